$wb = $excel.ActiveWorkbook

# --- Rename sheets ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16509960927734134"
$wb.Worksheets.Item(2).Name = "NB_TO-16509960945174143"
$wb.Worksheets.Item(3).Name = "RS_TO-16509960945174143"
$wb.Worksheets.Item(4).Name = "TOL_TO-1650996094565419"
$wb.Worksheets.Item(5).Name = "vSAT_TO-1650996094637413"

# --- Sheet 1 (GNG) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-1650996092741415.csv"
$ws1.Range("B3").Value = "GNG_stims-16509960927573867.csv"
$ws1.Range("B4").Value = "go_stims-16509960927573867.csv"
$ws1.Range("B5").Value = "GNG_stims-16509960927734134.csv"

# --- Sheet 2 (NB) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-16509960944933822.csv"
$ws2.Range("B3").Value = "ZB-match_9-16509960930694163.csv"
$ws2.Range("B4").Value = "TB-16509960940774133.csv"
$ws2.Range("B5").Value = "ZB-match_0-16509960931733959.csv"
$ws2.Range("B6").Value = "OB-16509960934133797.csv"
$ws2.Range("B7").Value = "ZB-match_3-16509960931413805.csv"
$ws2.Range("B8").Value = "TB-1650996094133387.csv"
$ws2.Range("B9").Value = "OB-16509960935174253.csv"
$ws2.Range("B10").Value = "OB-16509960933733811.csv"

# --- Sheet 3 (RS) --- no data changes

# --- Sheet 4 (TOL) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16509960945334184.csv"
$ws4.Range("B3").Value = "ZM_stims-16509960945174143.csv"
$ws4.Range("B4").Value = "MM_stims-16509960945493784.csv"
$ws4.Range("B5").Value = "ZM_stims-16509960945334184.csv"
$ws4.Range("B6").Value = "MM_stims-1650996094565419.csv"
$ws4.Range("B7").Value = "ZM_stims-16509960945493784.csv"

# --- Sheet 5 (vSAT) ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-1650996094565419.csv"
$ws5.Range("B3").Value = "SAT_stims-16509960945894132.csv"
$ws5.Range("B4").Value = "vSAT_stims-16509960946214201.csv"
$ws5.Range("B5").Value = "vSAT_stims-1650996094605414.csv"
